$d = $word.ActiveDocument

# Step 1: remove the old "_GoBack" bookmark (currently located near "772A NAS")
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Step 2: replace "day one, but will receive" with "day one, will receive"
$d.Content.Find.Execute(
    "fighting on day one, but will receive",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "fighting on day one, will receive",
    2)

# Step 3: insert the _GoBack bookmark right before "will receive" in the new text
$findRange = $d.Content
$findRange.Find.Execute("will receive her baptism of fire", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkRange = $d.Range($findRange.Start, $findRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
